$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Planilha1"

# Add new header in AE1 (shift existing "Metodo de pagamento" to AF1)
$ws.Range("AF1").Value = $ws.Range("AE1").Value2
$ws.Range("AE1").Value = "Data para envio"

# Copy the date-format style from AD1/AD2 onto AE1/AE2 (AE1 + AE2 both get style 1)
$ws.Range("AD1:AD2").Copy()
$ws.Range("AE1:AE2").PasteSpecial(-4122)

# Restore values (paste special format shouldn't clobber values, but just in case set again)
$ws.Range("AE1").Value = "Data para envio"

# Column width for AE (new "Data para envio" column), same stored width as AD (10.7265625)
$ws.Columns("AE").ColumnWidth = 9.893229166666666
